$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.782.78'
$ws.Range('E2').Value = '  +0.28%  '

# Row 3
$ws.Range('D3').Value = '2.570.12'
$ws.Range('E3').Value = '  +1.67%  '

# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').Value = '''302.67'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.43%  '

# Row 6
$ws.Range('D6').Value = '''97.31'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.28%  '

# Row 7
$ws.Range('D7').Value = '''0.574'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.97%  '

# Row 8
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('D9').Value = '''0.547'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.75%  '

# Row 10
$ws.Range('D10').Value = '''36.39'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.16%  '

# Row 11
$ws.Range('D11').Value = '''0.0809'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.86%  '

# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.116'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +8.82%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''7.58'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.60%  '

# Row 14
$ws.Range('D14').Value = '2.457.93'
$ws.Range('E14').Value = '  -2.67%  '

# Row 15
$ws.Range('D15').Value = '''0.883'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.85%  '

# Row 16
$ws.Range('D16').Value = '''14.40'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.70%  '

# Row 17
$ws.Range('D17').Value = '42.824.30'
$ws.Range('E17').Value = '  +0.37%  '

# Row 18
$ws.Range('B18').Value = 'InternetComputer(DFINITY)'
$ws.Range('C18').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D18').Value = '''13.32'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +8.56%  '

# Row 19
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0985'
$ws.Range('E19').Value = '  +2.83%  '

# Row 20
$ws.Range('D20').Value = '''6.63'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.83%  '

# Row 21
$ws.Range('D21').Value = '''71.72'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.12%  '

# Row 22
$ws.Range('D22').Value = '''258.03'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.50%  '

# Row 23
$ws.Range('D23').Value = '''2.95'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.78%  '

# Row 24
$ws.Range('E24').Value = '  -0.95%  '

# Row 25
$ws.Range('D25').Value = '''28.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.98%  '

# Row 26
$ws.Range('E26').Value = '  -0.02%  '

# Row 27
$ws.Range('B27').Value = 'InjectiveProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D27').Value = '''39.48'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +10.39%  '

# Row 28
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''10.11'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.61%  '

# Row 29
$ws.Range('E29').Value = '  -1.68%  '

# Row 30
$ws.Range('D30').Value = '''6.01'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.24%  '

# Row 31
$ws.Range('D31').Value = '''155.97'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.97%  '

# Row 32
$ws.Range('E32').Value = '  +1.86%  '

# Row 33
$ws.Range('E33').Value = '  +2.12%  '

# Row 34
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '''3.34'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.23%  '

# Row 35
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '''0.0803'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.95%  '

# Row 36
$ws.Range('D36').Value = '''18.17'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +12.98%  '

# Row 37
$ws.Range('E37').Value = '  +0.35%  '

# Row 38
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '''24.44'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.98%  '

# Row 39
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '''0.119'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.34%  '

# Row 40
$ws.Range('D40').Value = '''2.08'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +31.07%  '

# Row 41
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '''3.87'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.87%  '

# Row 42
$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').Value = '''3.38'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.22%  '

# Row 43
$ws.Range('D43').Value = '''0.0307'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.63%  '

# Row 44
$ws.Range('D44').Value = '2.091.67'
$ws.Range('E44').Value = '  +2.11%  '

# Row 45
$ws.Range('E45').Value = '  +0.08%  '

# Row 46
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').Value = '''88.61'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.56%  '

# Row 47
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '''9.28'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +6.19%  '

# Row 48
$ws.Range('D48').Value = '''77.62'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +12.38%  '

# Row 49
$ws.Range('D49').Value = '2.821.17'
$ws.Range('E49').Value = '  +1.84%  '

# Row 50
$ws.Range('D50').Value = '''104.46'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.77%  '

# Row 51
$ws.Range('D51').Value = '''0.190'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.16%  '

